# Generate Report for Handoff
# Refresh the handoff-generation timestamps and mark the handoff "priority"
# (ht = handoff type) for every localization file that was (re)packaged in
# this run, across the Overview summary sheet and each per-locale sheet.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 11, 12, 13, 14)

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-09-06 06:24:56"
}

# zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-09-06 06:24:51"
}

# de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-09-06 06:24:56"
}
